$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update IG version/date, publisher, and replace the
# duplicated "Contact" row with a single "Jurisdiction" row -----------------
$meta = $wb.Worksheets.Item("Metadata")

# The sheet currently has two identical "Contact" rows (10 and 11). Remove
# one of them entirely so the sheet shrinks from 21 to 20 data rows.
$meta.Rows(11).Delete()

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Publication date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a real value
$meta.Range("B9").Value = "Alvearie Team"

# The old "Contact" row becomes a "Jurisdiction" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Sheet "Elements": give the root Extension row a real short/definition -
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Employee Retirement Date"
$elements.Range("L2").Value = "Date of retirement for the employee or contract holder"
